$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 125.375
$ws.Range("I9").Value = 120.6
$ws.Range("J9").Value = 133.33333
$ws.Range("K9").Value = 120.6
$ws.Range("L9").Value = 133.33333
$ws.Range("M9").Value = 48.40000000000001
$ws.Range("N9").Value = -471.33333
$ws.Range("H12").Value = 850.36365
$ws.Range("J12").Value = 2774.5
$ws.Range("L12").Value = 2774.5
$ws.Range("N12").Value = -3114.5
$ws.Range("H15").Value = 3352.0222
$ws.Range("I15").Value = 3352.0222
$ws.Range("K15").Value = 10056.0666
$ws.Range("M15").Value = -9887.0666
$ws.Range("H32").Value = 7263.3335
$ws.Range("J32").Value = 8045.25
$ws.Range("L32").Value = 8045.25
$ws.Range("N32").Value = -8697.25
$ws.Range("H55").Value = 364.26315
$ws.Range("J55").Value = 406.7
$ws.Range("L55").Value = 406.7
$ws.Range("N55").Value = -834.7
$ws.Range("H64").Value = 7778.2
$ws.Range("J64").Value = 14666.667
$ws.Range("L64").Value = 14666.667
$ws.Range("N64").Value = -15162.667
$ws.Range("H67").Value = 7778.2
$ws.Range("J67").Value = 14666.667
$ws.Range("L67").Value = 14666.667
$ws.Range("N67").Value = -16382.667
$ws.Range("H80").Value = 8373.08
$ws.Range("I80").Value = 6883.467
$ws.Range("J80").Value = 10607.5
$ws.Range("K80").Value = 20650.401
$ws.Range("L80").Value = 31822.5
$ws.Range("M80").Value = -19652.401
$ws.Range("N80").Value = -33818.5
$ws.Range("H83").Value = 8373.08
$ws.Range("I83").Value = 6883.467
$ws.Range("J83").Value = 10607.5
$ws.Range("K83").Value = 61951.20299999999
$ws.Range("L83").Value = 95467.5
$ws.Range("M83").Value = -56959.20299999999
$ws.Range("N83").Value = -105451.5
$ws.Range("H86").Value = 1603.6
$ws.Range("I86").Value = 1833.3334
$ws.Range("J86").Value = 1259
$ws.Range("K86").Value = 1833.3334
$ws.Range("L86").Value = 1259
$ws.Range("M86").Value = -710.3334
$ws.Range("N86").Value = -3505
$ws.Range("H87").Value = 64985
$ws.Range("J87").Value = 64985
$ws.Range("L87").Value = 64985
$ws.Range("N87").Value = -67481
$ws.Range("H89").Value = 1603.6
$ws.Range("I89").Value = 1833.3334
$ws.Range("J89").Value = 1259
$ws.Range("K89").Value = 9166.666999999999
$ws.Range("L89").Value = 6295
$ws.Range("M89").Value = -3550.666999999999
$ws.Range("N89").Value = -17527
$ws.Range("H90").Value = 64985
$ws.Range("J90").Value = 64985
$ws.Range("L90").Value = 194955
$ws.Range("N90").Value = -207435
$ws.Range("H92").Value = 1007.0476
$ws.Range("I92").Value = 647.2778
$ws.Range("J92").Value = 3165.6667
$ws.Range("K92").Value = 647.2778
$ws.Range("L92").Value = 3165.6667
$ws.Range("M92").Value = 600.7222
$ws.Range("N92").Value = -5661.6667
$ws.Range("H98").Value = 1204.1177
$ws.Range("I98").Value = 1204.1177
$ws.Range("K98").Value = 1204.1177
$ws.Range("M98").Value = 293.8823
$ws.Range("H104").Value = 291.5
$ws.Range("I104").Value = 83
$ws.Range("J104").Value = 500
$ws.Range("K104").Value = 249
$ws.Range("L104").Value = 1500
$ws.Range("M104").Value = 1498
$ws.Range("N104").Value = -4994
$ws.Range("H122").Value = 1204.1177
$ws.Range("I122").Value = 1204.1177
$ws.Range("K122").Value = 3612.3531
$ws.Range("M122").Value = -1162.3531
$ws.Range("H134").Value = 175000
$ws.Range("J134").Value = 175000
$ws.Range("L134").Value = 175000
$ws.Range("N134").Value = -185140
$ws.Range("H138").Value = 6112.4683
$ws.Range("J138").Value = 4741.5
$ws.Range("L138").Value = 14224.5
$ws.Range("N138").Value = -24504.5
$ws.Range("H140").Value = 55000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 1913.7142
$ws.Range("J141").Value = 1498.5
$ws.Range("L141").Value = 4495.5
$ws.Range("N141").Value = -14855.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1580.2222
$ws.Range("I13").Value = 1177
$ws.Range("J13").Value = 2386.6667
$ws.Range("K13").Value = 1177
$ws.Range("L13").Value = 2386.6667
$ws.Range("M13").Value = -1033
$ws.Range("N13").Value = -2674.6667
$ws.Range("H30").Value = 1679.6
$ws.Range("I30").Value = 849.5
$ws.Range("K30").Value = 849.5
$ws.Range("M30").Value = -699.5
$ws.Range("H32").Value = 11237.445
$ws.Range("I32").Value = 9567.662
$ws.Range("J32").Value = 32666.334
$ws.Range("K32").Value = 9567.662
$ws.Range("L32").Value = 32666.334
$ws.Range("M32").Value = -9280.662
$ws.Range("N32").Value = -33240.334
$ws.Range("H45").Value = 98147.52
$ws.Range("I45").Value = 145146.36
$ws.Range("K45").Value = 145146.36
$ws.Range("M45").Value = -144769.36
$ws.Range("H74").Value = 3571.1667
$ws.Range("J74").Value = 6242.3335
$ws.Range("L74").Value = 6242.3335
$ws.Range("N74").Value = -7990.3335
$ws.Range("H77").Value = 3571.1667
$ws.Range("J77").Value = 6242.3335
$ws.Range("L77").Value = 31211.6675
$ws.Range("N77").Value = -39947.6675
$ws.Range("H97").Value = 925.0769
$ws.Range("I97").Value = 913.56525
$ws.Range("K97").Value = 913.56525
$ws.Range("M97").Value = -417.56525
$ws.Range("H102").Value = 4979.4736
$ws.Range("I102").Value = 2355.4614
$ws.Range("K102").Value = 2355.4614
$ws.Range("M102").Value = -733.4614000000001
$ws.Range("H122").Value = 1765.6666
$ws.Range("I122").Value = 1594.625
$ws.Range("K122").Value = 4783.875
$ws.Range("M122").Value = -2333.875
$ws.Range("H132").Value = 2228.3076
$ws.Range("I132").Value = 2247.3333
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6741.999899999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -4211.999899999999
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1484.5
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1484.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1484.5
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -1764.5
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4471
$ws.Range("H20").Value = 3482.5264
$ws.Range("I20").Value = 2500.7693
$ws.Range("K20").Value = 2500.7693
$ws.Range("M20").Value = -2253.7693
$ws.Range("H86").Value = 2008.4445
$ws.Range("I86").Value = 1766
$ws.Range("J86").Value = 2493.3333
$ws.Range("K86").Value = 1766
$ws.Range("L86").Value = 2493.3333
$ws.Range("M86").Value = -643
$ws.Range("N86").Value = -4739.3333
$ws.Range("H89").Value = 2008.4445
$ws.Range("I89").Value = 1766
$ws.Range("J89").Value = 2493.3333
$ws.Range("K89").Value = 8830
$ws.Range("L89").Value = 12466.6665
$ws.Range("M89").Value = -3214
$ws.Range("N89").Value = -23698.6665
$ws.Range("H94").Value = 251.66667
$ws.Range("I94").Value = 251.66667
$ws.Range("K94").Value = 251.66667
$ws.Range("M94").Value = 199.33333
$ws.Range("H134").Value = 2646.2058
$ws.Range("I134").Value = 2159.9033
$ws.Range("J134").Value = 7671.3335
$ws.Range("K134").Value = 6479.7099
$ws.Range("L134").Value = 23014.0005
$ws.Range("M134").Value = -3944.7099
$ws.Range("N134").Value = -28084.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 845.7778
$ws.Range("I22").Value = 344.7143
$ws.Range("J22").Value = 2599.5
$ws.Range("K22").Value = 344.7143
$ws.Range("L22").Value = 2599.5
$ws.Range("M22").Value = 5.28570000000002
$ws.Range("N22").Value = -3299.5
$ws.Range("H31").Value = 4849.5625
$ws.Range("I31").Value = 3061.75
$ws.Range("K31").Value = 3061.75
$ws.Range("M31").Value = -2766.75
$ws.Range("H34").Value = 4849.5625
$ws.Range("I34").Value = 3061.75
$ws.Range("K34").Value = 3061.75
$ws.Range("M34").Value = -2859.75
$ws.Range("H99").Value = 5227.1665
$ws.Range("I99").Value = 3783.0833
$ws.Range("J99").Value = 8115.3335
$ws.Range("K99").Value = 3783.0833
$ws.Range("L99").Value = 8115.3335
$ws.Range("M99").Value = -2285.0833
$ws.Range("N99").Value = -11111.3335
$ws.Range("H107").Value = 1174.1072
$ws.Range("I107").Value = 893.84
$ws.Range("J107").Value = 3509.6667
$ws.Range("K107").Value = 893.84
$ws.Range("L107").Value = 3509.6667
$ws.Range("M107").Value = 1026.16
$ws.Range("N107").Value = -7349.6667
$ws.Range("H126").Value = 5227.1665
$ws.Range("I126").Value = 3783.0833
$ws.Range("J126").Value = 8115.3335
$ws.Range("K126").Value = 11349.2499
$ws.Range("L126").Value = 24346.0005
$ws.Range("M126").Value = -8879.249899999999
$ws.Range("N126").Value = -29286.0005
$ws.Range("H134").Value = 2686.3333
$ws.Range("I134").Value = 1478.1143
$ws.Range("K134").Value = 4434.3429
$ws.Range("M134").Value = -1899.3429
$ws.Range("H140").Value = 517957.4
$ws.Range("J140").Value = 517957.4
$ws.Range("L140").Value = 517957.4
$ws.Range("N140").Value = -528317.4
$ws.Range("H141").Value = 214956.3
$ws.Range("J141").Value = 232173.67
$ws.Range("L141").Value = 232173.67
$ws.Range("N141").Value = -242533.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 524
$ws.Range("I23").Value = 500
$ws.Range("J23").Value = 530
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1590
$ws.Range("M23").Value = -1265
$ws.Range("N23").Value = -2060
$ws.Range("H34").Value = 583.6667
$ws.Range("J34").Value = 900
$ws.Range("L34").Value = 2700
$ws.Range("N34").Value = -2868
$ws.Range("H39").Value = 1462.375
$ws.Range("J39").Value = 1500
$ws.Range("L39").Value = 4500
$ws.Range("N39").Value = -5088
$ws.Range("H46").Value = 580
$ws.Range("J46").Value = 600
$ws.Range("L46").Value = 1800
$ws.Range("N46").Value = -1982
$ws.Range("H55").Value = 8194.111000000001
$ws.Range("J55").Value = 10250
$ws.Range("L55").Value = 30750
$ws.Range("N55").Value = -31104
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -314
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("K89").Value = 4500
$ws.Range("M89").Value = 1428
$ws.Range("H107").Value = 2024.3572
$ws.Range("I107").Value = 983.6667
$ws.Range("J107").Value = 2804.875
$ws.Range("K107").Value = 2951.0001
$ws.Range("L107").Value = 8414.625
$ws.Range("M107").Value = -1031.0001
$ws.Range("N107").Value = -12254.625
$ws.Range("H114").Value = 41668610
$ws.Range("I114").Value = 66667536
$ws.Range("J114").Value = 3733
$ws.Range("K114").Value = 200002608
$ws.Range("L114").Value = 11199
$ws.Range("M114").Value = -199999354
$ws.Range("N114").Value = -17707
$ws.Range("H122").Value = 5255.853
$ws.Range("J122").Value = 5639.8076
$ws.Range("L122").Value = 50758.2684
$ws.Range("N122").Value = -55658.2684
$ws.Range("H131").Value = 3354.0356
$ws.Range("I131").Value = 1228.7778
$ws.Range("J131").Value = 4360.737
$ws.Range("K131").Value = 3686.3334
$ws.Range("L131").Value = 13082.211
$ws.Range("M131").Value = 1353.6666
$ws.Range("N131").Value = -23162.211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1040.16
$ws.Range("J97").Value = 1110.2
$ws.Range("L97").Value = 1110.2
$ws.Range("N97").Value = -2102.2
$ws.Range("H122").Value = 3102.1765
$ws.Range("I122").Value = 3160.6924
$ws.Range("K122").Value = 9482.0772
$ws.Range("M122").Value = -7032.0772
$ws.Range("H132").Value = 5793.4443
$ws.Range("I132").Value = 4279
$ws.Range("J132").Value = 7005
$ws.Range("K132").Value = 12837
$ws.Range("L132").Value = 21015
$ws.Range("M132").Value = -10307
$ws.Range("N132").Value = -26075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6772.625
$ws.Range("I40").Value = 6377.6875
$ws.Range("K40").Value = 6377.6875
$ws.Range("M40").Value = -6241.6875
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H122").Value = 6635.5
$ws.Range("I122").Value = 6684.385
$ws.Range("K122").Value = 20053.155
$ws.Range("M122").Value = -17603.155
$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920
$ws.Range("H132").Value = 271149.75
$ws.Range("I132").Value = 360168.1
$ws.Range("J132").Value = 12996.5
$ws.Range("K132").Value = 1080504.3
$ws.Range("L132").Value = 38989.5
$ws.Range("M132").Value = -1077974.3
$ws.Range("N132").Value = -44049.5
$ws.Range("H136").Value = 117655304
$ws.Range("I136").Value = 62508764
$ws.Range("K136").Value = 187526292
$ws.Range("M136").Value = -187523742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7318
$ws.Range("J15").Value = 7326.6
$ws.Range("L15").Value = 7326.6
$ws.Range("N15").Value = -7902.6
$ws.Range("H46").Value = 76715.164
$ws.Range("J46").Value = 76715.164
$ws.Range("L46").Value = 76715.164
$ws.Range("N46").Value = -77177.164
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H52").Value = 13666.667
$ws.Range("I52").Value = 3000
$ws.Range("J52").Value = 19000
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 19000
$ws.Range("M52").Value = -2774
$ws.Range("N52").Value = -19452
$ws.Range("H62").Value = 11879.6
$ws.Range("I62").Value = 9800
$ws.Range("J62").Value = 12399.5
$ws.Range("K62").Value = 9800
$ws.Range("L62").Value = 12399.5
$ws.Range("M62").Value = -9176
$ws.Range("N62").Value = -13647.5
$ws.Range("H65").Value = 11879.6
$ws.Range("I65").Value = 9800
$ws.Range("J65").Value = 12399.5
$ws.Range("K65").Value = 49000
$ws.Range("L65").Value = 61997.5
$ws.Range("M65").Value = -45880
$ws.Range("N65").Value = -68237.5
$ws.Range("H107").Value = 1667.9524
$ws.Range("I107").Value = 1417.8334
$ws.Range("J107").Value = 2001.4445
$ws.Range("K107").Value = 4253.5002
$ws.Range("L107").Value = 6004.333500000001
$ws.Range("M107").Value = -2333.5002
$ws.Range("N107").Value = -9844.333500000001
$ws.Range("H122").Value = 4395.421
$ws.Range("J122").Value = 6240.6665
$ws.Range("L122").Value = 18721.9995
$ws.Range("N122").Value = -23621.9995
$ws.Range("H124").Value = 63768.832
$ws.Range("J124").Value = 63768.832
$ws.Range("L124").Value = 63768.832
$ws.Range("N124").Value = -73588.83199999999
$ws.Range("H126").Value = 5133.7896
$ws.Range("I126").Value = 4268.375
$ws.Range("J126").Value = 5763.1816
$ws.Range("K126").Value = 12805.125
$ws.Range("L126").Value = 17289.5448
$ws.Range("M126").Value = -10335.125
$ws.Range("N126").Value = -22229.5448
$ws.Range("H132").Value = 141589.39
$ws.Range("I132").Value = 166147.66
$ws.Range("K132").Value = 498442.98
$ws.Range("M132").Value = -495912.98
$ws.Range("H134").Value = 76715.164
$ws.Range("J134").Value = 76715.164
$ws.Range("L134").Value = 230145.492
$ws.Range("N134").Value = -235215.492
$ws.Range("H136").Value = 2177.6326
$ws.Range("I136").Value = 1404.6364
$ws.Range("K136").Value = 4213.9092
$ws.Range("M136").Value = -1663.9092
$ws.Range("H140").Value = 55947.25
$ws.Range("J140").Value = 55947.25
$ws.Range("L140").Value = 55947.25
$ws.Range("N140").Value = -66307.25
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
